$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.293.51"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.267.49"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'496.50"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "'128.99"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("D8").Value = "'0.527"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").Value = "'0.0954"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("E11").Value = "  +4.14%  "
$ws.Range("D12").Value = "'4.79"
$ws.Range("E12").Value = "  +3.74%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "'22.92"
$ws.Range("E13").Value = "  +5.25%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.667.03"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "54.337.27"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "2.256.59"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "'10.26"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").Value = "'301.94"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").Value = "'61.10"
$ws.Range("E23").Value = "  -2.27%  "
$ws.Range("D24").Value = "'0.995"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Value = "'7.34"
$ws.Range("E26").Value = "  +3.61%  "
$ws.Range("D27").Value = "'170.56"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'5.96"
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("D31").Value = "'1.09"
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").Value = "'17.83"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").Value = "'0.941"
$ws.Range("E35").Value = "  +8.94%  "
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").Value = "'3.71"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D41").Value = "'124.99"
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("D44").Value = "'0.0896"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").Value = "'0.548"
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("D46").Value = "'241.68"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "'16.14"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D51").Value = "'0.934"
$ws.Range("E51").Value = "  -0.56%  "
